$d = $word.ActiveDocument

# --- Part 1: split the title run "AbstractFactory Subject" into three runs
#     ("AbstractFactory", " ", "Subject"), each of the two words wrapped in
#     w:proofErr spellStart/spellEnd markers (as Word's spell checker would
#     emit after re-checking the text), keeping formatting identical. ---
$titlePara = $d.Paragraphs(1)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:color w:val="7030A0"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:color w:val="7030A0"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>AbstractFactory</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:color w:val="7030A0"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:color w:val="7030A0"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Subject</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$titlePara.Range.InsertXML($titleXml)

# --- Part 2: the final paragraph ("Comment peut-on modéliser ce problème ?")
#     moves its paragraph-mark run formatting (rPr) into pPr, and two new
#     paragraphs are appended after it: a "Réflexion : " heading paragraph
#     (same style) and a body paragraph with a tab + the new question. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tailXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="7030A0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="7030A0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Comment peut-on modéliser ce problème ?</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p>' + `
  '<w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="7030A0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="7030A0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Réflexion : </w:t></w:r>' + `
  '</w:p>' + `
  '<w:p>' + `
  '<w:r><w:rPr><w:color w:val="7030A0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Imaginons que nous prenons les données d’une base de données. Comment adapter le patron pour qu’il aille chercher les informations directement en base ?</w:t></w:r>' + `
  '</w:p>'
$lastPara.Range.InsertXML($tailXml)
